$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row/column of the existing data (header is row 1)
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# New columns go immediately to the right of the existing data (AD, AE, AF)
$winsCol = $lastCol + 1
$lossesCol = $lastCol + 2
$tiesCol = $lastCol + 3

# Copy the format of an existing header cell (A1) onto the new header cells
$ws.Range("A1").Copy()
$headerRange = $ws.Range($ws.Cells.Item(1, $winsCol), $ws.Cells.Item(1, $tiesCol))
$headerRange.PasteSpecial(-4122) # xlPasteFormats

# Set the header labels for the new columns
$ws.Cells.Item(1, $winsCol).Value = "Wins"
$ws.Cells.Item(1, $lossesCol).Value = "Losses"
$ws.Cells.Item(1, $tiesCol).Value = "Ties"

# Fill in the team record values for every data row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $winsCol).Value = 93
    $ws.Cells.Item($r, $lossesCol).Value = 69
    $ws.Cells.Item($r, $tiesCol).Value = 0
}
